# Fruta / hortaliza, semanal
# Insert a new weekly record at row 31 (pushing the existing rows 31-47 down
# to 32-48) for "Agrícola del Norte S.A. de Arica" - Pera - Winter Nelis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 31..47 down to 32..48, leaving a blank row 31.
$ws.Rows(31).Insert()

# Populate the new row 31 with the new weekly observation.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C31").Value = 'Arica y Parinacota'
$ws.Range("D31").Value = 44981
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 'Fruta'
$ws.Range("G31").Value = 100104
$ws.Range("H31").Value = 'Frutos de pepita'
$ws.Range("I31").Value = 100104005
$ws.Range("J31").Value = 'Pera'
$ws.Range("K31").Value = 'Winter Nelis'
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 26000
$ws.Range("P31").Value = 25500
$ws.Range("Q31").Value = '$/caja 20 kilos empedrada'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 1275
$ws.Range("T31").Value = 20
